# Update data: 10 June 2022
# Adds the latest month (date serial 44682 = 2022-05-01) to both the
# "Canada" sheet (national figure) and the "Province" sheet (one row per
# province/territory), matching the existing layout of the two tables.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Canada")
$ws2 = $wb.Worksheets.Item("Province")

$dateFmt = "d-mmm-yy"
$newDate = 44682

# ---------------------------------------------------------------------
# Sheet "Canada": single new row 30
# ---------------------------------------------------------------------
$r = 30

$ws1.Cells.Item($r, 1).Value = $newDate
$ws1.Cells.Item($r, 1).NumberFormat = $dateFmt

$ws1.Cells.Item($r, 2).Value = "Canada"
$ws1.Cells.Item($r, 2).NumberFormat = $dateFmt

$ws1.Cells.Item($r, 4).Value = 1057.8
$ws1.Cells.Item($r, 5).Value = 1093.6

$ws1.Cells.Item($r, 3).Formula = "=(D" + $r + "-E" + $r + ")/E" + $r + "*100"

# ---------------------------------------------------------------------
# Sheet "Province": ten new rows, 282-291 (one per province/territory,
# same order as every previous month's block)
# ---------------------------------------------------------------------
$provinces = @(
    @{ Name = "Newfoundland & Labrador"; D = 25.7;                 E = 32.9;                 Styled = $true  },
    @{ Name = "Prince Edward Island";    D = 7.2;                  E = 7.6;                  Styled = $false },
    @{ Name = "Nova Scotia";             D = 34.299999999999997;   E = 33;                   Styled = $false },
    @{ Name = "New Brunswick";           D = 27.9;                 E = 29.7;                 Styled = $false },
    @{ Name = "Quebec";                  D = 192.9;                E = 225.8;                Styled = $false },
    @{ Name = "Ontario";                 D = 447.7;                E = 410.8;                Styled = $false },
    @{ Name = "Manitoba";                D = 32.799999999999997;   E = 35.4;                 Styled = $false },
    @{ Name = "Saskatchewan";            D = 29;                   E = 32.299999999999997;   Styled = $false },
    @{ Name = "Alberta";                 D = 131.6;                E = 164.4;                Styled = $false },
    @{ Name = "British Columbia";        D = 128.6;                E = 121.7;                Styled = $false }
)

$startRow = 282
for ($i = 0; $i -lt $provinces.Count; $i++) {
    $row = $startRow + $i
    $p = $provinces[$i]

    $ws2.Cells.Item($row, 1).Value = $newDate
    $ws2.Cells.Item($row, 1).NumberFormat = $dateFmt

    $ws2.Cells.Item($row, 2).Value = $p.Name
    if ($p.Styled) {
        $ws2.Cells.Item($row, 2).NumberFormat = $dateFmt
    }

    $ws2.Cells.Item($row, 4).Value = $p.D
    $ws2.Cells.Item($row, 5).Value = $p.E

    $ws2.Cells.Item($row, 3).Formula = "=(D" + $row + "-E" + $row + ")/E" + $row + "*100"
}

# ---------------------------------------------------------------------
# Restore the selection / active-sheet state the workbook was saved with
# (sheet1 cursor on the new last row, sheet2 cursor one column past it
# on the new last row, "Province" remains the active/visible sheet).
# ---------------------------------------------------------------------
$ws1.Range("A30").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("D292").Select() | Out-Null
